# Removed bad decimals from excel file
# The sheet stores rows of: A=index, B=operand1, C=operator(text), D=operand2, E=result.
# A batch of rows had their operand / result values (and in a few cases the
# operator itself) corrected so the stored result is an exact, "clean"
# number instead of a rounded/truncated repeating decimal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: 15 / 15 = 1
$ws.Range("B7").Value = 15
$ws.Range("D7").Value = 15
$ws.Range("E7").Value = 1

# Row 12: 20 / 10 = 2
$ws.Range("B12").Value = 20
$ws.Range("D12").Value = 10
$ws.Range("E12").Value = 2

# Row 21: 15 / 2 = 7.5
$ws.Range("B21").Value = 15
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 7.5

# Row 25: 60 / 5 = 12
$ws.Range("B25").Value = 60
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 12

# Row 39: operator changed from "/" to "*" -> 5 * 11 = 55
$ws.Range("B39").Value = 5
$ws.Range("C39").Value = "*"
$ws.Range("D39").Value = 11
$ws.Range("E39").Value = 55

# Row 41: operator changed from "/" to "+" -> 14 + 9 = 23
$ws.Range("B41").Value = 14
$ws.Range("C41").Value = "+"
$ws.Range("D41").Value = 9
$ws.Range("E41").Value = 23

# Row 42: 16 / 8 = 2
$ws.Range("B42").Value = 16
$ws.Range("D42").Value = 8
$ws.Range("E42").Value = 2

# Row 48: operator changed from "/" to "-" -> 11 - 18 = -7
$ws.Range("B48").Value = 11
$ws.Range("C48").Value = "-"
$ws.Range("D48").Value = 18
$ws.Range("E48").Value = -7

# Row 56: 2 / 3 = 0.6667
$ws.Range("B56").Value = 2
$ws.Range("D56").Value = 3
$ws.Range("E56").Value = 0.6667

# Row 66: 12 * 12 = 144
$ws.Range("B66").Value = 12
$ws.Range("D66").Value = 12
$ws.Range("E66").Value = 144

# Row 69: 18 / 9 = 2
$ws.Range("B69").Value = 18
$ws.Range("D69").Value = 9
$ws.Range("E69").Value = 2

# Row 78: 12 / 6 = 2
$ws.Range("B78").Value = 12
$ws.Range("D78").Value = 6
$ws.Range("E78").Value = 2

# Row 85: 22 / 2 = 11
$ws.Range("B85").Value = 22
$ws.Range("D85").Value = 2
$ws.Range("E85").Value = 11

# Row 86: 32 / 4 = 8
$ws.Range("B86").Value = 32
$ws.Range("D86").Value = 4
$ws.Range("E86").Value = 8

# Row 99: 36 / 4 = 9
$ws.Range("B99").Value = 36
$ws.Range("D99").Value = 4
$ws.Range("E99").Value = 9

# Scroll the view down and move the selection, matching the editor's final
# cursor position when they finished cleaning up the sheet.
$ws.Range("A73").Select()
$ws.Range("E67").Select()
